# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Thu Sep  5 17:55:55 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.414.00"
$ws.Range("E2").Value = "  -1.95%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.378.19"
$ws.Range("E3").Value = "  -1.44%  "
# Row 4
$ws.Range("E4").Value = "  -0.39%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.16"
$ws.Range("E5").Value = "  -1.17%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.46"
$ws.Range("E6").Value = "  -1.71%  "
# Row 7
$ws.Range("E7").Value = "  +0.16%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  -2.31%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.386.43"
$ws.Range("E9").Value = "  -2.65%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0986"
$ws.Range("E10").Value = "  +0.14%  "
# Row 11
$ws.Range("E11").Value = "  +0.75%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  +1.21%  "
# Row 13
$ws.Range("E13").Value = "  +0.46%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.799.61"
$ws.Range("E14").Value = "  -1.69%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.368.71"
$ws.Range("E15").Value = "  -2.21%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.67"
$ws.Range("E16").Value = "  -1.30%  "
# Row 17
$ws.Range("E17").Value = "  -1.31%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.469.11"
$ws.Range("E18").Value = "  +0.08%  "
# Row 19
$ws.Range("E19").Value = "  -2.50%  "
# Row 20
$ws.Range("E20").Value = "  -2.60%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.60"
$ws.Range("E21").Value = "  -2.17%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("E22").Value = "  -1.49%  "
# Row 23
$ws.Range("E23").Value = "  +0.46%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.72"
$ws.Range("E24").Value = "  -0.98%  "
# Row 25
$ws.Range("E25").Value = "  +0.41%  "
# Row 26
$ws.Range("E26").Value = "  -3.70%  "
# Row 27
$ws.Range("E27").Value = "  -4.73%  "
# Row 28
$ws.Range("E28").Value = "  -3.99%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.49"
$ws.Range("E29").Value = "  -1.03%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0716"
$ws.Range("E30").Value = "  -2.97%  "
# Row 31
$ws.Range("E31").Value = "  -3.25%  "
# Row 32
$ws.Range("E32").Value = "  +0.23%  "
# Row 33
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.77"
$ws.Range("E33").Value = "  -7.15%  "
# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.09"
$ws.Range("E34").Value = "  -4.45%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  +0.52%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.58"
$ws.Range("E36").Value = "  -2.39%  "
# Row 37
$ws.Range("E37").Value = "  -5.82%  "
# Row 38
$ws.Range("E38").Value = "  -2.79%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.03"
$ws.Range("E39").Value = "  -1.29%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.795"
$ws.Range("E40").Value = "  -3.19%  "
# Row 41
$ws.Range("E41").Value = "  -4.60%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "131.07"
$ws.Range("E42").Value = "  -2.43%  "
# Row 43
$ws.Range("E43").Value = "  -1.63%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.77"
$ws.Range("E44").Value = "  -2.84%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.567"
$ws.Range("E45").Value = "  -0.88%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0906"
$ws.Range("E46").Value = "  -1.48%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "242.21"
$ws.Range("E47").Value = "  -6.63%  "
# Row 48
$ws.Range("E48").Value = "  -2.40%  "
# Row 49
$ws.Range("E49").Value = "  -2.63%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.15"
$ws.Range("E50").Value = "  -0.23%  "
# Row 51
$ws.Range("E51").Value = "  -2.96%  "

"done"